$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet
$ws.Name = "test_ProductsComparison"

# Update header row: B1 and C1 used to describe separate product numbers;
# both columns now simply represent "product"
$ws.Range("B1").Value = "product"
$ws.Range("C1").Value = "product"

$wb.Save()
